# Apply weekly fruit/vegetable price updates (row data shuffled per diff)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44400
$ws.Range("H2").Value = 'Española'
$ws.Range("J2").Value = 70
$ws.Range("K2").Value = 15000
$ws.Range("L2").Value = 15000
$ws.Range("M2").Value = 15000
$ws.Range("N2").Value = '$/caja 30 unidades'
$ws.Range("P2").Value = 500
$ws.Range("Q2").Value = 30

# Row 3
$ws.Range("D3").Value = 44446
$ws.Range("J3").Value = 120
$ws.Range("K3").Value = 16000
$ws.Range("L3").Value = 16000
$ws.Range("M3").Value = 16000
$ws.Range("P3").Value = 400

# Row 4
$ws.Range("D4").Value = 44421
$ws.Range("H4").Value = 'Española'
$ws.Range("K4").Value = 16500
$ws.Range("L4").Value = 16500
$ws.Range("M4").Value = 16500
$ws.Range("N4").Value = '$/caja 30 unidades'
$ws.Range("P4").Value = 550
$ws.Range("Q4").Value = 30

# Row 5
$ws.Range("D5").Value = 44495
$ws.Range("J5").Value = 130
$ws.Range("K5").Value = 11000
$ws.Range("L5").Value = 11000
$ws.Range("M5").Value = 11000
$ws.Range("P5").Value = 275

# Row 6
$ws.Range("D6").Value = 44481
$ws.Range("I6").Value = 'Segunda'
$ws.Range("J6").Value = 120
$ws.Range("N6").Value = '$/caja 50 unidades'
$ws.Range("P6").Value = 220
$ws.Range("Q6").Value = 50

# Row 7
$ws.Range("D7").Value = 44407
$ws.Range("H7").Value = 'Española'
$ws.Range("I7").Value = 'Primera'
$ws.Range("J7").Value = 100
$ws.Range("K7").Value = 18000
$ws.Range("L7").Value = 18000
$ws.Range("M7").Value = 18000
$ws.Range("N7").Value = '$/caja 30 unidades'
$ws.Range("P7").Value = 600
$ws.Range("Q7").Value = 30

# Row 8
$ws.Range("D8").Value = 44488
$ws.Range("H8").Value = 'Madrigal'
$ws.Range("J8").Value = 120
$ws.Range("K8").Value = 12000
$ws.Range("L8").Value = 12000
$ws.Range("M8").Value = 12000
$ws.Range("N8").Value = '$/caja 40 unidades'
$ws.Range("P8").Value = 300
$ws.Range("Q8").Value = 40

# Row 9
$ws.Range("D9").Value = 44176
$ws.Range("J9").Value = 80
$ws.Range("N9").Value = '$/caja 40 unidades'
$ws.Range("P9").Value = 275
$ws.Range("Q9").Value = 40

# Row 10
$ws.Range("D10").Value = 44390
$ws.Range("H10").Value = 'Española'
$ws.Range("J10").Value = 80
$ws.Range("K10").Value = 16000
$ws.Range("L10").Value = 16000
$ws.Range("M10").Value = 16000
$ws.Range("N10").Value = '$/caja 30 unidades'
$ws.Range("P10").Value = 533
$ws.Range("Q10").Value = 30

# Row 11
$ws.Range("D11").Value = 44484
$ws.Range("J11").Value = 110
$ws.Range("N11").Value = '$/caja 50 unidades'
$ws.Range("P11").Value = 220
$ws.Range("Q11").Value = 50

# Row 12
$ws.Range("D12").Value = 44386
$ws.Range("J12").Value = 30
$ws.Range("K12").Value = 15000
$ws.Range("L12").Value = 15000
$ws.Range("M12").Value = 15000
$ws.Range("P12").Value = 500

# Row 13
$ws.Range("D13").Value = 44491
$ws.Range("H13").Value = 'Madrigal'
$ws.Range("J13").Value = 200
$ws.Range("K13").Value = 11000
$ws.Range("L13").Value = 11000
$ws.Range("M13").Value = 11000
$ws.Range("N13").Value = '$/caja 40 unidades'
$ws.Range("P13").Value = 275
$ws.Range("Q13").Value = 40

# Row 14
$ws.Range("D14").Value = 44161
$ws.Range("H14").Value = 'Madrigal'
$ws.Range("J14").Value = 30
$ws.Range("K14").Value = 11000
$ws.Range("L14").Value = 11000
$ws.Range("M14").Value = 11000
$ws.Range("N14").Value = '$/caja 40 unidades'
$ws.Range("P14").Value = 275
$ws.Range("Q14").Value = 40

# Row 16
$ws.Range("D16").Value = 44418
$ws.Range("J16").Value = 80
$ws.Range("K16").Value = 16000
$ws.Range("L16").Value = 16000
$ws.Range("M16").Value = 16000
$ws.Range("P16").Value = 533

# Row 17
$ws.Range("D17").Value = 44166
$ws.Range("H17").Value = 'Madrigal'
$ws.Range("K17").Value = 10000
$ws.Range("L17").Value = 10000
$ws.Range("M17").Value = 10000
$ws.Range("N17").Value = '$/caja 40 unidades'
$ws.Range("P17").Value = 250
$ws.Range("Q17").Value = 40

Write-Host "Applied 111 cell updates"
